# Capstone_regPrecise/genomes.xlsx update
#
# The "Genomes" sheet lost two rows of scraped RegPrecise genome data
# (duplicates / entries removed upstream):
#   - genomeId 345 "Aeromonas salmonicida subsp. salmonicida A449" (row 12)
#   - genomeId 58  "Bacillus cereus ATCC 14579"                    (row 30)
#
# Deleting row 12 first shifts everything below it up by one, so the
# "Bacillus cereus" row (originally row 30) is now row 29 by the time we
# get to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(12).Delete()
$ws.Rows.Item(29).Delete()

# Leave the grid with row 29 (now "Bacillus clausii KSM-K16") selected,
# matching the author's final cursor position/selection.
$ws.Rows.Item(29).Select()
